$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
}

Set-TextValue "D2" "316.80"
Set-TextValue "E2" "2.09%"
Set-TextValue "G2" "22"
Set-TextValue "D3" "41.35"
Set-TextValue "E3" "1.19%"
Set-TextValue "G3" "22"
Set-TextValue "D4" "5.223"
Set-TextValue "E4" "2.01%"
Set-TextValue "G4" "22"
Set-TextValue "D5" "0.07637"
Set-TextValue "E5" "-0.35%"
Set-TextValue "G5" "22"
Set-TextValue "B6" "GateToken"
Set-TextValue "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "4.347"
Set-TextValue "E6" "1.39%"
Set-TextValue "G6" "22"
Set-TextValue "B7" "FTXToken"
Set-TextValue "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.663"
Set-TextValue "E7" "1.87%"
Set-TextValue "G7" "22"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9336"
Set-TextValue "E8" "1.90%"
Set-TextValue "G8" "22"
Set-TextValue "B9" "BTSEToken"
Set-TextValue "C9" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D9" "2.425"
Set-TextValue "E9" "-1.62%"
Set-TextValue "G9" "22"
Set-TextValue "B10" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D10" "0.1240"
Set-TextValue "E10" "-0.33%"
Set-TextValue "G10" "22"
Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1855"
Set-TextValue "E11" "2.73%"
Set-TextValue "G11" "22"
Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.09130"
Set-TextValue "E12" "0.38%"
Set-TextValue "G12" "22"
Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.04132"
Set-TextValue "E13" "-2.81%"
Set-TextValue "G13" "22"
Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.1054"
Set-TextValue "E14" "0.35%"
Set-TextValue "G14" "22"
Set-TextValue "B15" "BitForexToken"
Set-TextValue "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001276"
Set-TextValue "E15" "2.21%"
Set-TextValue "G15" "22"
Set-TextValue "B16" "TigerCash"
Set-TextValue "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.005960"
Set-TextValue "E16" "5.76%"
Set-TextValue "G16" "22"
Set-TextValue "B17" "UpBots"
Set-TextValue "C17" "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D17" "0.007491"
Set-TextValue "E17" "1,897.31%"
Set-TextValue "G17" "22"
Set-TextValue "B18" "HotbitToken"
Set-TextValue "C18" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D18" "0.004078"
Set-TextValue "E18" "-1.31%"
Set-TextValue "G18" "22"
Set-TextValue "B19" "LEO"
Set-TextValue "C19" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D19" "3.358"
Set-TextValue "E19" "-0.02%"
Set-TextValue "G19" "22"
Set-TextValue "B20" "BitpandaEcosystemToken"
Set-TextValue "C20" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D20" "0.3362"
Set-TextValue "E20" "1.45%"
Set-TextValue "G20" "22"
Set-TextValue "B21" "MCDex"
Set-TextValue "C21" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D21" "8.414"
Set-TextValue "E21" "21.64%"
Set-TextValue "G21" "22"
Set-TextValue "B22" "ProBitToken"
Set-TextValue "C22" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D22" "0.1361"
Set-TextValue "E22" "-2.23%"
Set-TextValue "G22" "22"
Set-TextValue "B23" "ZBToken"
Set-TextValue "C23" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D23" "0.2879"
Set-TextValue "E23" "5.41%"
Set-TextValue "G23" "22"
Set-TextValue "B24" "CoinExToken"
Set-TextValue "C24" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D24" "0.04044"
Set-TextValue "E24" "0.21%"
Set-TextValue "G24" "22"
Set-TextValue "B25" "BitKan"
Set-TextValue "C25" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D25" "0.001275"
Set-TextValue "E25" "0.57%"
Set-TextValue "G25" "22"
Set-TextValue "D26" "0.0001279"
Set-TextValue "E26" "0.98%"
Set-TextValue "G26" "22"
Set-TextValue "G27" "22"
Set-TextValue "G28" "22"
Set-TextValue "G29" "22"
Set-TextValue "G30" "22"
Set-TextValue "G31" "22"
Set-TextValue "G32" "22"
Set-TextValue "G33" "22"
Set-TextValue "G34" "22"
Set-TextValue "G35" "22"
Set-TextValue "G36" "22"
Set-TextValue "G37" "22"
Set-TextValue "D38" "0.02498"
Set-TextValue "E38" "2.71%"
Set-TextValue "G38" "22"
Set-TextValue "D39" "0.05241"
Set-TextValue "E39" "-0.31%"
Set-TextValue "G39" "22"
Set-TextValue "D40" "0.007816"
Set-TextValue "E40" "-0.16%"
Set-TextValue "G40" "22"
Set-TextValue "D41" "0.1298"
Set-TextValue "E41" "-0.86%"
Set-TextValue "G41" "22"
Set-TextValue "D42" "0.007093"
Set-TextValue "E42" "4.38%"
Set-TextValue "G42" "22"
Set-TextValue "D43" "0.002055"
Set-TextValue "E43" "11.84%"
Set-TextValue "G43" "22"
Set-TextValue "D44" "0.008234"
Set-TextValue "E44" "0.68%"
Set-TextValue "G44" "22"
Set-TextValue "D45" "0.3455"
Set-TextValue "E45" "3.22%"
Set-TextValue "G45" "22"
Set-TextValue "D46" "0.00006717"
Set-TextValue "E46" "-2.03%"
Set-TextValue "G46" "22"
Set-TextValue "E47" "0.82%"
Set-TextValue "G47" "22"
Set-TextValue "D48" "0.3845"
Set-TextValue "E48" "82.45%"
Set-TextValue "G48" "22"
Set-TextValue "D49" "0.004231"
Set-TextValue "E49" "3.29%"
Set-TextValue "G49" "22"
Set-TextValue "D50" "0.00002115"
Set-TextValue "E50" "0.82%"
Set-TextValue "G50" "22"
Set-TextValue "D51" "0.0002015"
Set-TextValue "E51" "0.82%"
Set-TextValue "G51" "22"
